# Actualización automática 2025-08-08 08:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("C13").Value = 129.6
$wsVentasPorGrupo.Range("C22").Value = "1 de 20"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F13").Value = 129.6
$wsVentaMensual.Range("F22").Value = 2918.78

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Target stored OOXML column width is 23; this engine's ColumnWidth setter
# adds a fixed 5/6 padding offset before writing the raw width, so we
# compensate to land exactly on 23.
$wsCumplimiento.Columns.Item(5).ColumnWidth = 22.166666666666668

$wsCumplimiento.Range("D2").Value = 129.6
$wsCumplimiento.Range("E2").Value = 652.8650105215589
$wsCumplimiento.Range("F2").Value = 0.1656304093567251

$wsCumplimiento.Range("D19").Value = 2918.78
$wsCumplimiento.Range("E19").Value = 56469.44762291769
$wsCumplimiento.Range("F19").Value = 0.04914745088088223
